# Add "Booking Time" (col G) and "Booking Status" (col H) fields to the
# reservation sheet, add a new reservation row (14), and normalize the
# previously-mixed-format rows 11-13 so they match the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Write H1 before G1 so the shared-string table gains "Booking Status"
# (index 10) ahead of "Booking Time" (index 11), matching the target file.
$ws.Range("H1").Value = "Booking Status"
$ws.Range("G1").Value = "Booking Time"

# --- New Booking Time column format ------------------------------------
$ws.Range("G2:G14").NumberFormat = "h:mm"

# --- Existing rows 2-10: just add the two new columns -------------------
$ws.Range("G2").Value = 0.45833333333333331
$ws.Range("H2").Value = "pending"

$ws.Range("G3").Value = 0.5
$ws.Range("H3").Value = "cancelled"

$ws.Range("G4").Value = 0.54166666666666696
$ws.Range("H4").Value = "confirmed"

$ws.Range("G5").Value = 0.58333333333333304
$ws.Range("H5").Value = "pending"

$ws.Range("G6").Value = 0.625
$ws.Range("H6").Value = "cancelled"

$ws.Range("G7").Value = 0.66666666666666596
$ws.Range("H7").Value = "confirmed"

$ws.Range("G8").Value = 0.70833333333333304
$ws.Range("H8").Value = "pending"

$ws.Range("G9").Value = 0.75
$ws.Range("H9").Value = "cancelled"

$ws.Range("G10").Value = 0.45833333333333331
$ws.Range("H10").Value = "confirmed"

# --- Rows 11-13: re-key A:F so the stray float/"n"-typed values loaded ---
# --- from the old file collapse back to plain integers, then add G:H ----
$ws.Range("C11:C13").NumberFormat = "yyyy\-mm\-dd"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = 45411
$ws.Range("D11").Value = 10
$ws.Range("E11").Value = "4 seat"
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 0.45833333333333331
$ws.Range("H11").Value = "pending"

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 10
$ws.Range("C12").Value = 45411
$ws.Range("D12").Value = 10
$ws.Range("E12").Value = "4 seat"
$ws.Range("F12").Value = 3
$ws.Range("G12").Value = 0.5
$ws.Range("H12").Value = "cancelled"

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 45412
$ws.Range("D13").Value = 4
$ws.Range("E13").Value = "2 seat"
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.54166666666666696
$ws.Range("H13").Value = "confirmed"

# --- Row 14: brand new reservation --------------------------------------
$ws.Range("C14").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 45413
$ws.Range("D14").Value = 4
$ws.Range("E14").Value = "2 seat"
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.58333333333333304
$ws.Range("H14").Value = "pending"

$ws.Range("K8").Select()
